$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "Ruins of Chernobog 'District 14'`n"
$ws.Range("C51").Value = "[name=`"FrostNova`"]  'They were your real parents.... and they died protecting you.'`n"
$ws.Range("C68").Value = "[name=`"FrostNova`"]  My grandmother was sentenced to the mines together as a 'harborer' of criminals, and moved with them.`n"
$ws.Range("C156").Value = "[name=`"FrostNova`"]  If you don't take Ursus at their word, you'll see that he's in good company. Although, there is only one that the Infected call 'Patriot'.`n"
$ws.Range("C161").Value = "[name=`"FrostNova`"]  However, during the 'Great Rebellion', his son was campaigning for the rights of the Infected. Long out of touch with his father, the son still sees him as an Ursus puppet.`n"
$ws.Range("C207").Value = "[name=`"FrostNova`"]  'Am I afraid of death?'`n"
$ws.Range("C220").Value = "[name=`"FrostNova`"]  'That's our big sister. She saved all our lives.'`n"
$ws.Range("C222").Value = "[name=`"FrostNova`"]  The children from the mine carried the Originium crystals I had made, and we became the 'Yeti Squad'.`n"
$ws.Range("C238").Value = "[name=`"FrostNova`"]  It's only because 'enemies kill each other' that we've fought as we have. That's all.`n"
$ws.Range("C242").Value = "[name=`"FrostNova`"]  'Reunion, first name. That. Infected, can rely on.'`n"
$ws.Range("C243").Value = "[name=`"FrostNova`"]  'No matter where, Reunion goes. No matter what, it does. We cannot, abandon it. To lose Reunion, is to lose, Infected conviction. Fight for: conviction.'`n"
$ws.Range("C252").Value = "[name=`"FrostNova`"]  After a dozen years of rolling around on the tundra, an Infected woman reached out to us. She said, 'Come with me, and let us break all chains.'`n"
$ws.Range("C255").Value = "[name=`"FrostNova`"]  But she said to us, 'Come with me.'`n"
$ws.Range("C273").Value = "[name=`"FrostNova`"]  'What are my chances of winning if I went up against her?'`n"
